$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the oldest sample row (row 22) - sliding window shifted by 6 new samples
$ws.Rows.Item(22).Delete()

# Update sensor reading columns (C:H) for rows 2-21 with the new sliding-window values
$ws.Range("C2").Value2 = 2.566667938232422
$ws.Range("D2").Value2 = -3.378203916549682
$ws.Range("E2").Value2 = 3.007539582252503
$ws.Range("F2").Value2 = -0.2000583708286285
$ws.Range("G2").Value2 = -0.1212567538022995
$ws.Range("H2").Value2 = -0.0207694191485643
$ws.Range("C3").Value2 = 3.106618106365205
$ws.Range("D3").Value2 = -3.249815458059311
$ws.Range("E3").Value2 = 3.031012719869614
$ws.Range("F3").Value2 = -0.1815796941518783
$ws.Range("G3").Value2 = -0.0572686158120632
$ws.Range("H3").Value2 = 0.08643743395805351
$ws.Range("C4").Value2 = 2.987140679359436
$ws.Range("D4").Value2 = -3.142817544937134
$ws.Range("E4").Value2 = 3.183629143238068
$ws.Range("F4").Value2 = -0.0739146918058395
$ws.Range("G4").Value2 = -0.1140790879726409
$ws.Range("H4").Value2 = 0.1067487001419067
$ws.Range("C5").Value2 = 2.434188187122345
$ws.Range("D5").Value2 = -3.181812554597855
$ws.Range("E5").Value2 = 3.162444919347763
$ws.Range("F5").Value2 = -0.0395535230636596
$ws.Range("G5").Value2 = -0.0899499058723449
$ws.Range("H5").Value2 = -0.0404698215425014
$ws.Range("C6").Value2 = 2.282221984863281
$ws.Range("D6").Value2 = -3.265003252029419
$ws.Range("E6").Value2 = 3.094355344772339
$ws.Range("F6").Value2 = -0.0148134818300604
$ws.Range("G6").Value2 = 0.1036943718791008
$ws.Range("H6").Value2 = -0.1157589629292488
$ws.Range("C7").Value2 = 2.110153055191039
$ws.Range("D7").Value2 = -3.195758980512619
$ws.Range("E7").Value2 = 3.138975620269776
$ws.Range("F7").Value2 = 0.5971207618713379
$ws.Range("G7").Value2 = 1.289536476135254
$ws.Range("H7").Value2 = -0.3637702465057373
$ws.Range("C8").Value2 = 1.555334329605102
$ws.Range("D8").Value2 = -2.938729083538055
$ws.Range("E8").Value2 = 3.47747951745987
$ws.Range("F8").Value2 = 1.519069194793701
$ws.Range("G8").Value2 = -0.4518875777721405
$ws.Range("H8").Value2 = -0.6734789609909058
$ws.Range("C9").Value2 = 0.6493126988410929
$ws.Range("D9").Value2 = -2.875420850515366
$ws.Range("E9").Value2 = 3.464587104320525
$ws.Range("F9").Value2 = 0.2113593816757202
$ws.Range("G9").Value2 = -0.3769038617610931
$ws.Range("H9").Value2 = 0.4825835525989532
$ws.Range("C10").Value2 = -0.6712930202484151
$ws.Range("D10").Value2 = -3.392556905746461
$ws.Range("E10").Value2 = 2.365111112594603
$ws.Range("F10").Value2 = 0.2393064647912979
$ws.Range("G10").Value2 = -0.8791878223419189
$ws.Range("H10").Value2 = -0.1872301995754242
$ws.Range("C11").Value2 = 0.5537151455879301
$ws.Range("D11").Value2 = -4.361428594589236
$ws.Range("E11").Value2 = 3.347476267814645
$ws.Range("F11").Value2 = -0.1922698318958282
$ws.Range("G11").Value2 = -0.9285151958465576
$ws.Range("H11").Value2 = 0.8594874143600464
$ws.Range("C12").Value2 = 2.452674245834349
$ws.Range("D12").Value2 = -4.409869003295896
$ws.Range("E12").Value2 = 5.084140586853025
$ws.Range("F12").Value2 = -3.570354700088501
$ws.Range("G12").Value2 = -0.7802276611328125
$ws.Range("H12").Value2 = -4.989242076873779
$ws.Range("C13").Value2 = 5.681596696376824
$ws.Range("D13").Value2 = -3.525700151920317
$ws.Range("E13").Value2 = 2.85166837722062
$ws.Range("F13").Value2 = -1.221577763557434
$ws.Range("G13").Value2 = 2.375196695327759
$ws.Range("H13").Value2 = -2.503631114959717
$ws.Range("C14").Value2 = 1.737989616393984
$ws.Range("D14").Value2 = -2.796510410308835
$ws.Range("E14").Value2 = 1.393881118297584
$ws.Range("F14").Value2 = 2.165364503860474
$ws.Range("G14").Value2 = 0.5566509366035461
$ws.Range("H14").Value2 = -0.4453207552433014
$ws.Range("C15").Value2 = -5.58917605876923
$ws.Range("D15").Value2 = -7.849099040031435
$ws.Range("E15").Value2 = 6.400659620761871
$ws.Range("F15").Value2 = 0.2727513313293457
$ws.Range("G15").Value2 = 0.5925393104553223
$ws.Range("H15").Value2 = 0.4948008358478546
$ws.Range("C16").Value2 = 2.166972637176496
$ws.Range("D16").Value2 = -11.98566874265667
$ws.Range("E16").Value2 = 9.183138275146462
$ws.Range("F16").Value2 = -0.2121229618787765
$ws.Range("G16").Value2 = 1.80510675907135
$ws.Range("H16").Value2 = 1.96942949295044
$ws.Range("C17").Value2 = -0.1177038192748863
$ws.Range("D17").Value2 = -2.770210593938835
$ws.Range("E17").Value2 = 3.860614097118379
$ws.Range("F17").Value2 = 0.2755002379417419
$ws.Range("G17").Value2 = 1.588096976280212
$ws.Range("H17").Value2 = 2.037540912628174
$ws.Range("C18").Value2 = 3.127950906753536
$ws.Range("D18").Value2 = -4.12096252441406
$ws.Range("E18").Value2 = 3.36216964721679
$ws.Range("F18").Value2 = 0.8869763612747192
$ws.Range("G18").Value2 = 0.8231409192085266
$ws.Range("H18").Value2 = 1.362993121147156
$ws.Range("C19").Value2 = 2.318384975194932
$ws.Range("D19").Value2 = -3.449181020259855
$ws.Range("E19").Value2 = 1.645497059822083
$ws.Range("F19").Value2 = 0.2379320114850998
$ws.Range("G19").Value2 = -0.7533495426177979
$ws.Range("H19").Value2 = 0.1786780804395675
$ws.Range("C20").Value2 = 3.076034724712372
$ws.Range("D20").Value2 = -2.935223340988159
$ws.Range("E20").Value2 = 1.488467574119568
$ws.Range("F20").Value2 = 0.6151412725448608
$ws.Range("G20").Value2 = 1.230893492698669
$ws.Range("H20").Value2 = -0.3686571717262268
$ws.Range("C21").Value2 = 4.307219874858854
$ws.Range("D21").Value2 = -2.422872281074523
$ws.Range("E21").Value2 = 1.113696080446245
$ws.Range("F21").Value2 = 0.1403462886810302
$ws.Range("G21").Value2 = 0.7915286421775818
$ws.Range("H21").Value2 = 0.00137444678694
